$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Octubre de 2020 a las 20:19"

# Update country names where the ranking swapped (row order stays the same,
# but two adjacent pairs exchanged places because their case counts crossed)
$ws.Range("A65").Value = "Libano"
$ws.Range("A66").Value = "Argelia"
$ws.Range("A72").Value = "Irlanda"
$ws.Range("A73").Value = "Azerbaiyan"

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# Update the statistics (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) for the refreshed data pull
Set-Row 4   7972747 24457 5113389 2639798 0 190 219560
Set-Row 5   7079426 27883 6104199 866654  0 202 108573
Set-Row 13  734974  16101 100828 601463  0 46  32683
Set-Row 25  325039  1586  273500 41840   0 8   9699
Set-Row 33  152404  2563  127407 22392   0 33  2605
Set-Row 65  53568   1010  23501  29608   0 4   459
Set-Row 66  53072   132   37170  14101   0 6   1801
Set-Row 72  42528   814   23364  17338   0 2   1826
Set-Row 73  41982   230   39357  2016    0 1   609
Set-Row 105 10894   35    9742   1117    0 1   35
Set-Row 147 3286    13    2527   627     0 0   132
